$wb = $excel.ActiveWorkbook

# Rename the four year/type analysis sheets to their new short English names.
$wsModified = $wb.Worksheets.Item("AÑO MODIFICACIÓN")
$wsModified.Name = "modified"

$wsPublished = $wb.Worksheets.Item("AÑO PUBLICACIÓN")
$wsPublished.Name = "published"

$wsType = $wb.Worksheets.Item("TIPO DE OBJETO STIX 2.1 ")
$wsType.Name = "type"

$wsCreated = $wb.Worksheets.Item("AÑO CREACION")
$wsCreated.Name = "created"

# Each of those sheets hosts one chart whose series still point at the old
# (now stale) quoted Spanish sheet names - repoint them at the new names
# while leaving the cached values / series titles untouched.
$chModified = $wsModified.ChartObjects().Item(1).Chart
$serModified = $chModified.SeriesCollection().Item(1)
$serModified.XValues = "=modified!`$B`$12:`$B`$13"
$serModified.Values = "=modified!`$C`$12:`$C`$13"

$chPublished = $wsPublished.ChartObjects().Item(1).Chart
$serPublished = $chPublished.SeriesCollection().Item(1)
$serPublished.XValues = "=published!`$B`$12:`$B`$13"
$serPublished.Values = "=published!`$C`$12:`$C`$13"

$chType = $wsType.ChartObjects().Item(1).Chart
$serType = $chType.SeriesCollection().Item(1)
$serType.XValues = "=type!`$B`$12:`$B`$13"
$serType.Values = "=type!`$C`$12:`$C`$13"

$chCreated = $wsCreated.ChartObjects().Item(1).Chart
$serCreated = $chCreated.SeriesCollection().Item(1)
$serCreated.XValues = "=created!`$B`$12:`$B`$14"
$serCreated.Values = "=created!`$D`$12:`$D`$14"

# Move the active/selected tab from "type" (previously 3rd tab, index 2) to
# "created" (now the 4th tab, index 3).
$wsCreated.Activate()
$wsCreated.Select()
